$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: new boiling point correlation data row (columns A:DQ)
$row14Values = @(0.8894, 0.9225, 0.6033, 0.2878, 1.7827, 1.8433, 1.7117, 1.7957, 1.8881, 3.1243, 0.9297, 1.6254, 1.9669, 1.9478, 1.7444, 3.2152, 4.4014, 3.5668, 3.8967, 2.8526, 3.636, 3.3953, 3.1459, 2.2536, 1.6249, 1.1557, 2.5892, 3.1656, 2.5983, 3.1376, 2.6127, 1.578, 2.1647, 1.2171, 5.4736, 6.28, 5.9234, 5.0525, 5.8337, 2.9637, 2.6948, 2.2073, 3.93, 3.56, 4.5797, 2.6293, 5.7619, 5.0767, 6.0837, 3.2914, 3.665, 2.6495, 2.3678, 2.5645, 1.7824, 0.9442, 7.2644, 1.288, 0.6115, 1.1739, 2.6446, 2.8881, 2.3086, 1.9163, 1.0081, 10.3428, 0.0, 0.0, 7.6904, 0.0, 6.7822, 5.5566, 5.4248, 3.6796, 3.6763, 2.6812, 5.7093, 5.826, -0.1157, -0.0489, 0.1798, 0.3189, 0.7273, 0.4745, 0.3563, 0.1919, 0.1957, 0.3489, 0.1589, 0.0668, -0.1406, -0.09, 0.0511, 0.6884, -0.1074, 0.0224, 0.092, 0.558, 0.0735, -0.1552, 0.7801, -0.2383, 0.4456, -0.1977, 0.0835, -0.5385, -0.6331, 1.4108, -0.069, 1.0682, 0.4247, 0.2499, 0.1134, -0.2596, 0.4408, -0.1168, -0.3201, -0.4453, -0.6776, -0.3678, 0.0)
for ($i = 0; $i -lt $row14Values.Length; $i++) {
    $cell = $ws.Cells.Item(14, $i + 1)
    $cell.Value = $row14Values[$i]
    $cell.Font.Color = 0
}

# New label row 31: "tb"
$ws.Range("A31").Value = "tb"

# Update sheet view: selection moves to A31, remove frozen/top-left scroll position
$ws.Range("A31").Select()
